$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance "O" for 1주차 (row2) and 3주차 (row4) in column C,
# matching the existing style used for 7주차/8주차 (C7/C8).
$ws.Range("C2").Value = "O"
$ws.Range("C4").Value = "O"

# Update the active selection to C6.
$ws.Range("C6").Select()
